$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Job adverts by profession" row to the refreshed SOC based job adverts data.
$ws.Range("A13").Value = "Job adverts by occupation"
$ws.Range("B13").Value = "<a href='https://www.ons.gov.uk/employmentandlabourmarket/peopleinwork/employmentandemployeetypes/datasets/labourdemandvolumesbystandardoccupationclassificationsoc2020uk'>ONS Textkernel</a>"
$ws.Range("C13").Value = "May 2023 (15/03/24)"

# Update selected cell to match the saved selection state in the source workbook.
$ws.Range("C14").Select()
